$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cell F1 with same style as E1 (bold header with border)
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("F1").Value = "time_taken"

# Fill F2:F181 with time_taken values (plain text, default style)
$arr = New-Object 'object[,]' 180,1
$arr[0,0] = "2021-10-05 10:50:10.144670"
$arr[1,0] = "2021-10-05 10:50:10.144681"
$arr[2,0] = "2021-10-05 10:50:10.144685"
$arr[3,0] = "2021-10-05 10:50:10.144687"
$arr[4,0] = "2021-10-05 10:50:10.144690"
$arr[5,0] = "2021-10-05 10:50:10.144693"
$arr[6,0] = "2021-10-05 10:50:10.144695"
$arr[7,0] = "2021-10-05 10:50:10.144698"
$arr[8,0] = "2021-10-05 10:50:10.144701"
$arr[9,0] = "2021-10-05 10:50:10.144703"
$arr[10,0] = "2021-10-05 10:50:10.144706"
$arr[11,0] = "2021-10-05 10:50:10.144708"
$arr[12,0] = "2021-10-05 10:50:10.144711"
$arr[13,0] = "2021-10-05 10:50:10.144713"
$arr[14,0] = "2021-10-05 10:50:10.144716"
$arr[15,0] = "2021-10-05 10:50:10.144718"
$arr[16,0] = "2021-10-05 10:50:10.144721"
$arr[17,0] = "2021-10-05 10:50:10.144724"
$arr[18,0] = "2021-10-05 10:50:10.144726"
$arr[19,0] = "2021-10-05 10:50:10.144728"
$arr[20,0] = "2021-10-05 10:50:10.144731"
$arr[21,0] = "2021-10-05 10:50:10.144733"
$arr[22,0] = "2021-10-05 10:50:10.144736"
$arr[23,0] = "2021-10-05 10:50:10.144738"
$arr[24,0] = "2021-10-05 10:50:10.144741"
$arr[25,0] = "2021-10-05 10:50:10.144744"
$arr[26,0] = "2021-10-05 10:50:10.144746"
$arr[27,0] = "2021-10-05 10:50:10.144749"
$arr[28,0] = "2021-10-05 10:50:10.144751"
$arr[29,0] = "2021-10-05 10:50:10.144754"
$arr[30,0] = "2021-10-05 10:50:10.144757"
$arr[31,0] = "2021-10-05 10:50:10.144759"
$arr[32,0] = "2021-10-05 10:50:10.144762"
$arr[33,0] = "2021-10-05 10:50:10.144764"
$arr[34,0] = "2021-10-05 10:50:10.144767"
$arr[35,0] = "2021-10-05 10:50:10.144769"
$arr[36,0] = "2021-10-05 10:50:10.144772"
$arr[37,0] = "2021-10-05 10:50:10.144774"
$arr[38,0] = "2021-10-05 10:50:10.144777"
$arr[39,0] = "2021-10-05 10:50:10.144779"
$arr[40,0] = "2021-10-05 10:50:10.144782"
$arr[41,0] = "2021-10-05 10:50:10.144785"
$arr[42,0] = "2021-10-05 10:50:10.144787"
$arr[43,0] = "2021-10-05 10:50:10.144790"
$arr[44,0] = "2021-10-05 10:50:10.144792"
$arr[45,0] = "2021-10-05 10:50:10.144795"
$arr[46,0] = "2021-10-05 10:50:10.144797"
$arr[47,0] = "2021-10-05 10:50:10.144800"
$arr[48,0] = "2021-10-05 10:50:10.144802"
$arr[49,0] = "2021-10-05 10:50:10.144804"
$arr[50,0] = "2021-10-05 10:50:10.144807"
$arr[51,0] = "2021-10-05 10:50:10.144809"
$arr[52,0] = "2021-10-05 10:50:10.144812"
$arr[53,0] = "2021-10-05 10:50:10.144815"
$arr[54,0] = "2021-10-05 10:50:10.144818"
$arr[55,0] = "2021-10-05 10:50:10.144820"
$arr[56,0] = "2021-10-05 10:50:10.144822"
$arr[57,0] = "2021-10-05 10:50:10.144825"
$arr[58,0] = "2021-10-05 10:50:10.144827"
$arr[59,0] = "2021-10-05 10:50:10.144830"
$arr[60,0] = "2021-10-05 10:50:10.144832"
$arr[61,0] = "2021-10-05 10:50:10.144835"
$arr[62,0] = "2021-10-05 10:50:10.144837"
$arr[63,0] = "2021-10-05 10:50:10.144840"
$arr[64,0] = "2021-10-05 10:50:10.144843"
$arr[65,0] = "2021-10-05 10:50:10.144846"
$arr[66,0] = "2021-10-05 10:50:10.144848"
$arr[67,0] = "2021-10-05 10:50:10.144851"
$arr[68,0] = "2021-10-05 10:50:10.144853"
$arr[69,0] = "2021-10-05 10:50:10.144856"
$arr[70,0] = "2021-10-05 10:50:10.144858"
$arr[71,0] = "2021-10-05 10:50:10.144861"
$arr[72,0] = "2021-10-05 10:50:10.144863"
$arr[73,0] = "2021-10-05 10:50:10.144866"
$arr[74,0] = "2021-10-05 10:50:10.144868"
$arr[75,0] = "2021-10-05 10:50:10.144870"
$arr[76,0] = "2021-10-05 10:50:10.144875"
$arr[77,0] = "2021-10-05 10:50:10.144878"
$arr[78,0] = "2021-10-05 10:50:10.144880"
$arr[79,0] = "2021-10-05 10:50:10.144883"
$arr[80,0] = "2021-10-05 10:50:10.144885"
$arr[81,0] = "2021-10-05 10:50:10.144888"
$arr[82,0] = "2021-10-05 10:50:10.144890"
$arr[83,0] = "2021-10-05 10:50:10.144893"
$arr[84,0] = "2021-10-05 10:50:10.144895"
$arr[85,0] = "2021-10-05 10:50:10.144898"
$arr[86,0] = "2021-10-05 10:50:10.144900"
$arr[87,0] = "2021-10-05 10:50:10.144903"
$arr[88,0] = "2021-10-05 10:50:10.144905"
$arr[89,0] = "2021-10-05 10:50:10.144908"
$arr[90,0] = "2021-10-05 10:50:10.144910"
$arr[91,0] = "2021-10-05 10:50:10.144913"
$arr[92,0] = "2021-10-05 10:50:10.144916"
$arr[93,0] = "2021-10-05 10:50:10.144919"
$arr[94,0] = "2021-10-05 10:50:10.144922"
$arr[95,0] = "2021-10-05 10:50:10.144925"
$arr[96,0] = "2021-10-05 10:50:10.144927"
$arr[97,0] = "2021-10-05 10:50:10.144930"
$arr[98,0] = "2021-10-05 10:50:10.144932"
$arr[99,0] = "2021-10-05 10:50:10.144935"
$arr[100,0] = "2021-10-05 10:50:10.144937"
$arr[101,0] = "2021-10-05 10:50:10.144940"
$arr[102,0] = "2021-10-05 10:50:10.144942"
$arr[103,0] = "2021-10-05 10:50:10.144945"
$arr[104,0] = "2021-10-05 10:50:10.144947"
$arr[105,0] = "2021-10-05 10:50:10.144950"
$arr[106,0] = "2021-10-05 10:50:10.144952"
$arr[107,0] = "2021-10-05 10:50:10.144955"
$arr[108,0] = "2021-10-05 10:50:10.144959"
$arr[109,0] = "2021-10-05 10:50:10.144962"
$arr[110,0] = "2021-10-05 10:50:10.144965"
$arr[111,0] = "2021-10-05 10:50:10.144967"
$arr[112,0] = "2021-10-05 10:50:10.144969"
$arr[113,0] = "2021-10-05 10:50:10.144972"
$arr[114,0] = "2021-10-05 10:50:10.144974"
$arr[115,0] = "2021-10-05 10:50:10.144977"
$arr[116,0] = "2021-10-05 10:50:10.144979"
$arr[117,0] = "2021-10-05 10:50:10.144982"
$arr[118,0] = "2021-10-05 10:50:10.144984"
$arr[119,0] = "2021-10-05 10:50:10.144987"
$arr[120,0] = "2021-10-05 10:50:10.144989"
$arr[121,0] = "2021-10-05 10:50:10.144992"
$arr[122,0] = "2021-10-05 10:50:10.144994"
$arr[123,0] = "2021-10-05 10:50:10.144996"
$arr[124,0] = "2021-10-05 10:50:10.144999"
$arr[125,0] = "2021-10-05 10:50:10.145001"
$arr[126,0] = "2021-10-05 10:50:10.145004"
$arr[127,0] = "2021-10-05 10:50:10.145006"
$arr[128,0] = "2021-10-05 10:50:10.145011"
$arr[129,0] = "2021-10-05 10:50:10.145014"
$arr[130,0] = "2021-10-05 10:50:10.145016"
$arr[131,0] = "2021-10-05 10:50:10.145019"
$arr[132,0] = "2021-10-05 10:50:10.145021"
$arr[133,0] = "2021-10-05 10:50:10.145023"
$arr[134,0] = "2021-10-05 10:50:10.145026"
$arr[135,0] = "2021-10-05 10:50:10.145028"
$arr[136,0] = "2021-10-05 10:50:10.145031"
$arr[137,0] = "2021-10-05 10:50:10.145033"
$arr[138,0] = "2021-10-05 10:50:10.145036"
$arr[139,0] = "2021-10-05 10:50:10.145038"
$arr[140,0] = "2021-10-05 10:50:10.145041"
$arr[141,0] = "2021-10-05 10:50:10.145043"
$arr[142,0] = "2021-10-05 10:50:10.145046"
$arr[143,0] = "2021-10-05 10:50:10.145048"
$arr[144,0] = "2021-10-05 10:50:10.145051"
$arr[145,0] = "2021-10-05 10:50:10.145053"
$arr[146,0] = "2021-10-05 10:50:10.145055"
$arr[147,0] = "2021-10-05 10:50:10.145058"
$arr[148,0] = "2021-10-05 10:50:10.145061"
$arr[149,0] = "2021-10-05 10:50:10.145063"
$arr[150,0] = "2021-10-05 10:50:10.145066"
$arr[151,0] = "2021-10-05 10:50:10.145068"
$arr[152,0] = "2021-10-05 10:50:10.145071"
$arr[153,0] = "2021-10-05 10:50:10.145073"
$arr[154,0] = "2021-10-05 10:50:10.145076"
$arr[155,0] = "2021-10-05 10:50:10.145078"
$arr[156,0] = "2021-10-05 10:50:10.145081"
$arr[157,0] = "2021-10-05 10:50:10.145083"
$arr[158,0] = "2021-10-05 10:50:10.145086"
$arr[159,0] = "2021-10-05 10:50:10.145088"
$arr[160,0] = "2021-10-05 10:50:10.145091"
$arr[161,0] = "2021-10-05 10:50:10.145093"
$arr[162,0] = "2021-10-05 10:50:10.145095"
$arr[163,0] = "2021-10-05 10:50:10.145098"
$arr[164,0] = "2021-10-05 10:50:10.145100"
$arr[165,0] = "2021-10-05 10:50:10.145103"
$arr[166,0] = "2021-10-05 10:50:10.145105"
$arr[167,0] = "2021-10-05 10:50:10.145108"
$arr[168,0] = "2021-10-05 10:50:10.145110"
$arr[169,0] = "2021-10-05 10:50:10.145113"
$arr[170,0] = "2021-10-05 10:50:10.145115"
$arr[171,0] = "2021-10-05 10:50:10.145118"
$arr[172,0] = "2021-10-05 10:50:10.145121"
$arr[173,0] = "2021-10-05 10:50:10.145124"
$arr[174,0] = "2021-10-05 10:50:10.145127"
$arr[175,0] = "2021-10-05 10:50:10.145130"
$arr[176,0] = "2021-10-05 10:50:10.145132"
$arr[177,0] = "2021-10-05 10:50:10.145135"
$arr[178,0] = "2021-10-05 10:50:10.145137"
$arr[179,0] = "2021-10-05 10:50:10.145140"
$ws.Range("F2:F181").Value = $arr
